$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a brand-new row at position 583. This shifts the previous rows
# 583..668 down to 584..669 (Excel's standard Insert-row behaviour), which
# is exactly the shift pattern shown across the whole diff (every row N's
# old data becomes row N+1's data).
$ws.Rows.Item(583).Insert()

# Populate the newly-inserted row 583 with the new weekly data point.
$ws.Range("A583").Value = 6
$ws.Range("B583").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C583").Value = "Metropolitana"
$ws.Range("D583").Value = 45077
$ws.Range("E583").Value = 13
$ws.Range("F583").Value = "Fruta"
$ws.Range("G583").Value = 100101
$ws.Range("H583").Value = "Berries"
$ws.Range("I583").Value = 100101001
$ws.Range("J583").Value = "Arándano (blue)"
$ws.Range("K583").Value = "Sin especificar"
$ws.Range("L583").Value = "Primera"
$ws.Range("M583").Value = 75
$ws.Range("N583").Value = 12000
$ws.Range("O583").Value = 12000
$ws.Range("P583").Value = 12000
$ws.Range("Q583").Value = "$/bandeja 2 kilos"
$ws.Range("R583").Value = "Provincia de Curicó"
$ws.Range("S583").Value = 6000
$ws.Range("T583").Value = 2
